# Fix typos / stray characters in the "Spatial Structure" and
# "Recurrent Mutation" option columns of the raw dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated closing parenthesis.
$ws.Cells.Replace("Continuous Space (1D or 2D))", "Continuous Space (1D or 2D)")
$ws.Cells.Replace("Metapopulation (Finite or Infinite))", "Metapopulation (Finite or Infinite)")

# Normalize "1D, 2D" -> "1D or 2D" for consistency with the other options.
$ws.Cells.Replace("Stepping Stone (1D, 2D)", "Stepping Stone (1D or 2D)")

# Fix spelling/casing typo: "adpative dynamics" -> "Adaptive Dynamics".
$ws.Cells.Replace("adpative dynamics", "Adaptive Dynamics")
